$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value to a cell without Excel coercing
# numeric-looking strings (e.g. "313.98") into a Number. Number-format the
# cell as Text for the write, then restore its original style so no visible
# formatting / style index changes are left behind.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" '27.294.64'
$ws.Range("E2").Value = '  +0.77%  '
Set-TextValue "D3" '1.850.11'
$ws.Range("E3").Value = '  +0.96%  '
Set-TextValue "D4" '1.001'
$ws.Range("E4").Value = '  -0.51%  '
Set-TextValue "D5" '313.98'
$ws.Range("E5").Value = '  +0.69%  '
Set-TextValue "D6" '1.001'
$ws.Range("E6").Value = '  -0.61%  '
Set-TextValue "D7" '0.4603'
$ws.Range("E7").Value = '  -1.13%  '
Set-TextValue "D8" '0.3707'
$ws.Range("E8").Value = '  -0.05%  '
Set-TextValue "D9" '0.07288'
$ws.Range("E9").Value = '  -1.15%  '
Set-TextValue "D10" '0.8838'
$ws.Range("E10").Value = '  +1.01%  '
Set-TextValue "D11" '19.89'
$ws.Range("E11").Value = '  -0.49%  '
Set-TextValue "D12" '0.07793'
$ws.Range("E12").Value = '  -0.73%  '
Set-TextValue "D13" '1.846.49'
$ws.Range("E13").Value = '  +3.86%  '
Set-TextValue "D14" '5.366'
$ws.Range("E14").Value = '  +0.20%  '
Set-TextValue "D15" '6.543'
$ws.Range("E15").Value = '  -0.48%  '
Set-TextValue "D16" '91.46'
$ws.Range("E16").Value = '  -0.52%  '
Set-TextValue "D17" '1.002'
$ws.Range("E17").Value = '  -0.46%  '
Set-TextValue "D18" '0.000008957'
$ws.Range("E18").Value = '  +1.11%  '
Set-TextValue "D19" '1.001'
$ws.Range("E19").Value = '  -0.37%  '
Set-TextValue "D20" '14.75'
$ws.Range("E20").Value = '  +0.71%  '
Set-TextValue "D21" '27.317.44'
$ws.Range("E21").Value = '  +2.30%  '
Set-TextValue "D22" '5.116'
$ws.Range("E22").Value = '  -0.46%  '
Set-TextValue "D23" '10.51'
$ws.Range("E23").Value = '  -0.64%  '
Set-TextValue "D24" '2.074.79'
$ws.Range("E24").Value = '  +7.22%  '
Set-TextValue "D25" '1.916'
$ws.Range("E25").Value = '  +4.93%  '
Set-TextValue "D26" '151.47'
$ws.Range("E26").Value = '  -0.59%  '
Set-TextValue "D27" '18.38'
$ws.Range("E27").Value = '  +0.28%  '
Set-TextValue "D28" '2.057'
$ws.Range("E28").Value = '  -2.00%  '
Set-TextValue "D29" '116.05'
$ws.Range("E29").Value = '  +0.61%  '
Set-TextValue "D30" '5.087'
$ws.Range("E30").Value = '  -0.22%  '
Set-TextValue "D31" '0.08828'
$ws.Range("E31").Value = '  -0.62%  '
Set-TextValue "D32" '3.127'
$ws.Range("E32").Value = '  +5.50%  '
Set-TextValue "D33" '0.7715'
$ws.Range("E33").Value = '  +5.86%  '
Set-TextValue "D34" '1.169'
$ws.Range("E34").Value = '  +2.63%  '
Set-TextValue "D35" '4.491'
$ws.Range("E35").Value = '  +1.05%  '
Set-TextValue "D36" '2.662'
$ws.Range("E36").Value = '  +6.43%  '
Set-TextValue "D37" '1.078'
$ws.Range("E37").Value = '  +0.33%  '
Set-TextValue "D38" '0.01956'
$ws.Range("E38").Value = '  +0.23%  '
Set-TextValue "D39" '0.05225'
$ws.Range("E39").Value = '  +0.00%  '
Set-TextValue "D40" '2.954'
$ws.Range("E40").Value = '  +0.70%  '
Set-TextValue "D41" '6.992'
$ws.Range("E41").Value = '  -1.99%  '
Set-TextValue "D42" '0.5134'
$ws.Range("E42").Value = '  -1.32%  '
Set-TextValue "D43" '0.1632'
$ws.Range("E43").Value = '  +0.09%  '
Set-TextValue "D44" '8.403'
$ws.Range("E44").Value = '  +1.94%  '
Set-TextValue "D45" '0.4804'
$ws.Range("E45").Value = '  -0.49%  '
Set-TextValue "D46" '10.27'
$ws.Range("E46").Value = '  +1.18%  '
Set-TextValue "D47" '1.000'
$ws.Range("E47").Value = '  -0.76%  '
Set-TextValue "D48" '102.52'
$ws.Range("E48").Value = '  +0.25%  '
Set-TextValue "D49" '1.648'
$ws.Range("E49").Value = '  +1.11%  '
Set-TextValue "D50" '0.06216'
$ws.Range("E50").Value = '  +0.06%  '
Set-TextValue "D51" '65.36'
$ws.Range("E51").Value = '  +2.17%  '
